$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Scratch helpers ---
# M1: holds a copy of the "MEDIUM THREAT - MONITOR" (yellow) style, captured
#     from J10 before its content/format gets overwritten below.
# M2: reformatted as Text ("@") so that date-like strings (e.g. "13-JAN-26")
#     assigned through it keep their literal text instead of being parsed into
#     a real Excel date serial number when pasted onward.
$mediumStyle = $ws.Range("M1")
$ws.Range("J10").Copy()
$mediumStyle.PasteSpecial(-4122)

$textScratch = $ws.Range("M2")
$textScratch.NumberFormat = "@"

function Set-LiteralText($cellRef, [string]$text) {
    $textScratch.Value = $text
    $textScratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

function Set-MediumThreatStyle($cellRef) {
    $mediumStyle.Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

function Set-LowThreatStyle($cellRef) {
    $ws.Range("J2").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

# --- Row data (after the edit) ---

# Row 2
Set-LiteralText "A2" '13-JAN-26'
$ws.Range("B2").Value = 'SM-433'
$ws.Range("C2").Value = 'Nile Air NP-131'
$ws.Range("D2").Value = 11341
$ws.Range("E2").Value = 11858
$ws.Range("F2").Value = -517
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 'LOW THREAT'
$ws.Range("K2").Value = 'EGP'

# Row 3
Set-LiteralText "A3" '13-JAN-26'
$ws.Range("B3").Value = 'SM-433'
$ws.Range("C3").Value = 'Nesma Airlines NE-150'
$ws.Range("D3").Value = 11618
$ws.Range("E3").Value = 11858
$ws.Range("F3").Value = -240
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 'LOW THREAT'
$ws.Range("K3").Value = 'EGP'

# Row 4
Set-LiteralText "A4" '13-JAN-26'
$ws.Range("B4").Value = 'SM-433'
$ws.Range("C4").Value = 'flynas XY-855'
$ws.Range("D4").Value = 13345
$ws.Range("E4").Value = 11858
$ws.Range("F4").Value = 1487
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = -10
$ws.Range("J4").Value = 'LOW THREAT'
$ws.Range("K4").Value = 'EGP'

# Row 5
Set-LiteralText "A5" '15-JAN-26'
$ws.Range("B5").Value = 'SM-433'
$ws.Range("C5").Value = 'Nesma Airlines NE-154'
$ws.Range("D5").Value = 11618
$ws.Range("E5").Value = 13155
$ws.Range("F5").Value = -1537
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 'LOW THREAT'
$ws.Range("K5").Value = 'EGP'

# Row 6
Set-LiteralText "A6" '15-JAN-26'
$ws.Range("B6").Value = 'SM-433'
$ws.Range("C6").Value = 'EgyptAir MS-681'
$ws.Range("D6").Value = 14693
$ws.Range("E6").Value = 13155
$ws.Range("F6").Value = 1538
$ws.Range("G6").Value = 46
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = -16
$ws.Range("J6").Value = 'LOW THREAT'
$ws.Range("K6").Value = 'EGP'

# Row 7
Set-LiteralText "A7" '16-JAN-26'
$ws.Range("B7").Value = 'SM-433'
$ws.Range("C7").Value = 'Nile Air NP-133'
$ws.Range("D7").Value = 11341
$ws.Range("E7").Value = 14453
$ws.Range("F7").Value = -3112
$ws.Range("G7").Value = 30
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = 0
Set-MediumThreatStyle "J7"
$ws.Range("J7").Value = 'MEDIUM THREAT - MONITOR'
$ws.Range("K7").Value = 'EGP'

# Row 8
Set-LiteralText "A8" '17-JAN-26'
$ws.Range("B8").Value = 'SM-433'
$ws.Range("C8").Value = 'Nesma Airlines NE-154'
$ws.Range("D8").Value = 11618
$ws.Range("E8").Value = 14453
$ws.Range("F8").Value = -2835
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = 30
$ws.Range("I8").Value = 0
Set-MediumThreatStyle "J8"
$ws.Range("J8").Value = 'MEDIUM THREAT - MONITOR'
$ws.Range("K8").Value = 'EGP'

# Row 9
Set-LiteralText "A9" '17-JAN-26'
$ws.Range("B9").Value = 'SM-433'
$ws.Range("C9").Value = 'Nesma Airlines NE-152'
$ws.Range("D9").Value = 12929
$ws.Range("E9").Value = 14453
$ws.Range("F9").Value = -1524
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 'LOW THREAT'
$ws.Range("K9").Value = 'EGP'

# Row 10
Set-LiteralText "A10" '22-JAN-26'
$ws.Range("B10").Value = 'SM-433'
$ws.Range("C10").Value = 'Nile Air NP-131'
$ws.Range("D10").Value = 10207
$ws.Range("E10").Value = 10887
$ws.Range("F10").Value = -680
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = 0
Set-LowThreatStyle "J10"
$ws.Range("J10").Value = 'LOW THREAT'
$ws.Range("K10").Value = 'EGP'

# Row 11
Set-LiteralText "A11" '22-JAN-26'
$ws.Range("B11").Value = 'SM-433'
$ws.Range("C11").Value = 'EgyptAir MS-681'
$ws.Range("D11").Value = 12450
$ws.Range("E11").Value = 10887
$ws.Range("F11").Value = 1563
$ws.Range("G11").Value = 46
$ws.Range("H11").Value = 30
$ws.Range("I11").Value = -16
Set-LowThreatStyle "J11"
$ws.Range("J11").Value = 'LOW THREAT'
$ws.Range("K11").Value = 'EGP'

# Row 12
Set-LiteralText "A12" '24-JAN-26'
$ws.Range("B12").Value = 'SM-433'
$ws.Range("C12").Value = 'Air Arabia Egypt E5-315'
$ws.Range("D12").Value = 9119
$ws.Range("E12").Value = 9904
$ws.Range("F12").Value = -785
$ws.Range("G12").Value = 30
$ws.Range("H12").Value = 30
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 'LOW THREAT'
$ws.Range("K12").Value = 'EGP'

# Row 13
Set-LiteralText "A13" '27-JAN-26'
$ws.Range("B13").Value = 'SM-433'
$ws.Range("C13").Value = 'EgyptAir MS-681'
$ws.Range("D13").Value = 8316
$ws.Range("E13").Value = 6792
$ws.Range("F13").Value = 1524
$ws.Range("G13").Value = 46
$ws.Range("H13").Value = 30
$ws.Range("I13").Value = -16
$ws.Range("J13").Value = 'LOW THREAT'
$ws.Range("K13").Value = 'EGP'

# Row 14
Set-LiteralText "A14" '30-JAN-26'
$ws.Range("B14").Value = 'SM-433'
$ws.Range("C14").Value = 'EgyptAir MS-681'
$ws.Range("D14").Value = 8417
$ws.Range("E14").Value = 6792
$ws.Range("F14").Value = 1625
$ws.Range("G14").Value = 46
$ws.Range("H14").Value = 30
$ws.Range("I14").Value = -16
$ws.Range("J14").Value = 'LOW THREAT'
$ws.Range("K14").Value = 'EGP'

# Row 15
Set-LiteralText "A15" '31-JAN-26'
$ws.Range("B15").Value = 'SM-433'
$ws.Range("C15").Value = 'Air Arabia Egypt E5-315'
$ws.Range("D15").Value = 6293
$ws.Range("E15").Value = 7434
$ws.Range("F15").Value = -1141
$ws.Range("G15").Value = 30
$ws.Range("H15").Value = 30
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 'LOW THREAT'
$ws.Range("K15").Value = 'EGP'

# Row 16
Set-LiteralText "A16" '31-JAN-26'
$ws.Range("B16").Value = 'SM-433'
$ws.Range("C16").Value = 'Nesma Airlines NE-154'
$ws.Range("D16").Value = 6830
$ws.Range("E16").Value = 7434
$ws.Range("F16").Value = -604
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 30
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 'LOW THREAT'
$ws.Range("K16").Value = 'EGP'

# Row 17
Set-LiteralText "A17" '31-JAN-26'
$ws.Range("B17").Value = 'SM-433'
$ws.Range("C17").Value = 'EgyptAir MS-681'
$ws.Range("D17").Value = 8316
$ws.Range("E17").Value = 7434
$ws.Range("F17").Value = 882
$ws.Range("G17").Value = 46
$ws.Range("H17").Value = 30
$ws.Range("I17").Value = -16
$ws.Range("J17").Value = 'LOW THREAT'
$ws.Range("K17").Value = 'EGP'

# --- Clean up scratch cells, then drop the now-obsolete trailing rows ---
$textScratch.Clear()
$mediumStyle.Clear()
$ws.Rows("18:24").Delete()
